$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 2..N) is sorted ascending by column A (ID) to
# produce a better-sampled ordering of the selected papers.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$dataRange = $ws.Range("A2:G" + $lastRow)
$keyRange = $ws.Range("A2:A" + $lastRow)

$dataRange.Sort($keyRange)
